# Updated cryptos list (GitHub Actions scheduled refresh)
# Applies the latest price / 1h volume figures, plus a ranking swap
# (Cronos -> BabyDogeCoin) in row 51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces a text value so Excel does not reinterpret
# decimal-looking prices (e.g. "0.515") as numbers and strip
# trailing zeros; entries that already contain two dots (e.g.
# "29.128.40") can never be parsed as numbers so no prefix is needed.
$ws.Range("D2").Value = "29.128.40"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "1.578.25"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'211.94"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'0.515"
$ws.Range("E6").Value = "  +6.25%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "'25.87"
$ws.Range("E8").Value = "  +9.29%  "
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "1.804.45"
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("D13").Value = "1.560.03"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "29.153.38"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "'62.30"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D18").Value = "'238.20"
$ws.Range("E18").Value = "  +4.91%  "
$ws.Range("D19").Value = "'7.45"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("E20").Value = "  +2.78%  "
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("E23").Value = "  +4.12%  "
$ws.Range("E24").Value = "  +5.24%  "
$ws.Range("D25").Value = "'153.34"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("E26").Value = "  +4.35%  "
$ws.Range("D27").Value = "'15.11"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "'0.0465"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "1.423.06"
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").Value = "'2.74"
$ws.Range("E37").Value = "  +6.00%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").Value = "'0.528"
$ws.Range("E40").Value = "  +3.54%  "
$ws.Range("E41").Value = "  +2.70%  "
$ws.Range("D42").Value = "'53.13"
$ws.Range("E42").Value = "  +25.76%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").Value = "'0.788"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").Value = "'0.0471"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").Value = "'64.41"
$ws.Range("E46").Value = "  +4.25%  "
$ws.Range("D47").Value = "'5.35"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "1.716.22"
$ws.Range("D49").Value = "'0.849"
$ws.Range("E49").Value = "  -6.36%  "
$ws.Range("D50").Value = "'85.76"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0101"
$ws.Range("E51").Value = "  -1.73%  "
